## updated template and example for importing samples
##
## The workbook's "Scats" sheet is the sample-upload template/example.
## Column D used to be "genotype_id" (left blank for every sample row);
## it is repurposed as "sample_type", and every existing example row is
## back-filled with the sample type "excrement" (all rows in this example
## sheet are scat/excrement samples).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column D header from "genotype_id" to "sample_type".
$ws.Range("D1").Value = "sample_type"

# Fill the new sample_type column for every data row (2-25) with
# "excrement", matching the plain left/top-aligned "General" text style
# already used by the other populated text columns (e.g. column A).
$sampleTypeRange = $ws.Range("D2:D25")
$sampleTypeRange.Value = "excrement"
$sampleTypeRange.NumberFormat = "General"
$sampleTypeRange.HorizontalAlignment = -4131
$sampleTypeRange.VerticalAlignment = -4160
$sampleTypeRange.WrapText = $false

# Leave the freshly-filled column selected (and scrolled back to the top
# of the sheet), matching the state the sheet was left in after the edit.
$ws.Range("A1").Select()
$sampleTypeRange.Select()
